$wb = $excel.ActiveWorkbook

# --- Cell value updates ---
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("A2").Value = "ffff3f313e35-bb77-428d-8718-d48193819a0b.md"
$ws.Range("B2").Value = "e2e\ffff3f313e35-bb77-428d-8718-d48193819a0b.md"
$ws.Range("G2").Value = "2016-08-31 09:18:37"
$ws.Range("A3").Value = "ffffff23b151d0-47d5-4d0f-aef1-d540bdbaff4a.md"
$ws.Range("B3").Value = "e2e\ffffff23b151d0-47d5-4d0f-aef1-d540bdbaff4a.md"
$ws.Range("A4").Value = "78e8bd2f-a385-4f19-88eb-6ad4215eb7ce.md"
$ws.Range("B4").Value = "e2e\78e8bd2f-a385-4f19-88eb-6ad4215eb7ce.md"
$ws.Range("E4").Value = "Ready for handoff"
$ws.Range("F4").Value = "Ready for handoff"
$ws.Range("G4").Value = "2016-08-31 09:20:56"

$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("A2").Value = "ffff3f313e35-bb77-428d-8718-d48193819a0b.md"
$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("G2").Value = "2047bb2e-acca-4a3c-8dfd-f8a610485fda.fb05da26239bb986374d99f05efe0bc1d19a2bf4.zh-cn.xlf"
$ws.Range("H2").Value = "2016-08-31 09:18:32"
$ws.Range("I2").Value = "2047bb2e-acca-4a3c-8dfd-f8a610485fda.md"
$ws.Range("J2").Value = "2047bb2e-acca-4a3c-8dfd-f8a610485fda.fb05da26239bb986374d99f05efe0bc1d19a2bf4.zh-cn.xlf"
$ws.Range("K2").Value = "2016-08-31 09:18:49"
$ws.Range("A3").Value = "ffffff23b151d0-47d5-4d0f-aef1-d540bdbaff4a.md"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("F3").Value = "'True"
$ws.Range("A4").Value = "78e8bd2f-a385-4f19-88eb-6ad4215eb7ce.md"
$ws.Range("C4").Value = "Ready for handoff"
$ws.Range("F4").Value = "'False"
$ws.Range("G4").Value = "78e8bd2f-a385-4f19-88eb-6ad4215eb7ce.7072c78345e460eca7557c4583f70bbec535734c.zh-cn.xlf"
$ws.Range("H4").Value = "2016-08-31 09:20:51"
$ws.Range("I4").Value = "78e8bd2f-a385-4f19-88eb-6ad4215eb7ce.md"
$ws.Range("J4").Value = "78e8bd2f-a385-4f19-88eb-6ad4215eb7ce.7072c78345e460eca7557c4583f70bbec535734c.zh-cn.xlf"
$ws.Range("K4").Value = "2016-08-31 09:20:03"
$ws.Range("P4").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5952b6da77cfabfe6add37ddd290d5d20161214c/e2e/78e8bd2f-a385-4f19-88eb-6ad4215eb7ce.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ec2147e48cc81f604574fa1f2ea3f511c0909e10/e2e/78e8bd2f-a385-4f19-88eb-6ad4215eb7ce.md."

$ws = $wb.Worksheets.Item("de-de")
$ws.Range("A2").Value = "ffff3f313e35-bb77-428d-8718-d48193819a0b.md"
$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("G2").Value = "2047bb2e-acca-4a3c-8dfd-f8a610485fda.fb05da26239bb986374d99f05efe0bc1d19a2bf4.de-de.xlf"
$ws.Range("H2").Value = "2016-08-31 09:18:37"
$ws.Range("I2").Value = "2047bb2e-acca-4a3c-8dfd-f8a610485fda.md"
$ws.Range("J2").Value = "2047bb2e-acca-4a3c-8dfd-f8a610485fda.fb05da26239bb986374d99f05efe0bc1d19a2bf4.de-de.xlf"
$ws.Range("K2").Value = "2016-08-31 09:18:56"
$ws.Range("A3").Value = "ffffff23b151d0-47d5-4d0f-aef1-d540bdbaff4a.md"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("F3").Value = "'True"
$ws.Range("A4").Value = "78e8bd2f-a385-4f19-88eb-6ad4215eb7ce.md"
$ws.Range("C4").Value = "Ready for handoff"
$ws.Range("F4").Value = "'False"
$ws.Range("G4").Value = "78e8bd2f-a385-4f19-88eb-6ad4215eb7ce.7072c78345e460eca7557c4583f70bbec535734c.de-de.xlf"
$ws.Range("H4").Value = "2016-08-31 09:20:56"
$ws.Range("I4").Value = "78e8bd2f-a385-4f19-88eb-6ad4215eb7ce.md"
$ws.Range("J4").Value = "78e8bd2f-a385-4f19-88eb-6ad4215eb7ce.7072c78345e460eca7557c4583f70bbec535734c.de-de.xlf"
$ws.Range("K4").Value = "2016-08-31 09:20:33"
$ws.Range("P4").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5952b6da77cfabfe6add37ddd290d5d20161214c/e2e/78e8bd2f-a385-4f19-88eb-6ad4215eb7ce.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ec2147e48cc81f604574fa1f2ea3f511c0909e10/e2e/78e8bd2f-a385-4f19-88eb-6ad4215eb7ce.md."

# --- Hyperlink display text updates ---
$ws = $wb.Worksheets.Item("Overview")
foreach ($hl in $ws.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$B$2') {
        $hl.TextToDisplay = "e2e\ffff3f313e35-bb77-428d-8718-d48193819a0b.md"
    }
    elseif ($addr -eq '$B$3') {
        $hl.TextToDisplay = "e2e\ffffff23b151d0-47d5-4d0f-aef1-d540bdbaff4a.md"
    }
    elseif ($addr -eq '$B$4') {
        $hl.TextToDisplay = "e2e\78e8bd2f-a385-4f19-88eb-6ad4215eb7ce.md"
    }
}

$ws = $wb.Worksheets.Item("zh-cn")
foreach ($hl in $ws.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = "ffff3f313e35-bb77-428d-8718-d48193819a0b.md"
    }
    elseif ($addr -eq '$I$2') {
        $hl.TextToDisplay = "2047bb2e-acca-4a3c-8dfd-f8a610485fda.md"
    }
    elseif ($addr -eq '$A$3') {
        $hl.TextToDisplay = "ffffff23b151d0-47d5-4d0f-aef1-d540bdbaff4a.md"
    }
    elseif ($addr -eq '$I$3') {
        $hl.TextToDisplay = "2047bb2e-acca-4a3c-8dfd-f8a610485fda.md"
    }
    elseif ($addr -eq '$A$4') {
        $hl.TextToDisplay = "78e8bd2f-a385-4f19-88eb-6ad4215eb7ce.md"
    }
    elseif ($addr -eq '$I$4') {
        $hl.TextToDisplay = "78e8bd2f-a385-4f19-88eb-6ad4215eb7ce.md"
    }
}

$ws = $wb.Worksheets.Item("de-de")
foreach ($hl in $ws.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = "ffff3f313e35-bb77-428d-8718-d48193819a0b.md"
    }
    elseif ($addr -eq '$I$2') {
        $hl.TextToDisplay = "2047bb2e-acca-4a3c-8dfd-f8a610485fda.md"
    }
    elseif ($addr -eq '$A$3') {
        $hl.TextToDisplay = "ffffff23b151d0-47d5-4d0f-aef1-d540bdbaff4a.md"
    }
    elseif ($addr -eq '$I$3') {
        $hl.TextToDisplay = "2047bb2e-acca-4a3c-8dfd-f8a610485fda.md"
    }
    elseif ($addr -eq '$A$4') {
        $hl.TextToDisplay = "78e8bd2f-a385-4f19-88eb-6ad4215eb7ce.md"
    }
    elseif ($addr -eq '$I$4') {
        $hl.TextToDisplay = "78e8bd2f-a385-4f19-88eb-6ad4215eb7ce.md"
    }
}

# --- Column width updates ---
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Columns.Item(16).ColumnWidth = 40
$ws = $wb.Worksheets.Item("de-de")
$ws.Columns.Item(16).ColumnWidth = 40
